# Generate Report for Handoff
# - Update status text "In Translation" -> "Ready for handoff"
# - Bump the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# - Widen the Status-related columns to fit the new "Ready for handoff" text

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status values: "In Translation" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Timestamps ---
$overview.Range("G2").Value = "2016-09-07 09:31:39"
$dede.Range("H2").Value     = "2016-09-07 09:31:39"
$zhcn.Range("H2").Value     = "2016-09-07 09:31:34"

# --- Column widths (to fit the new, longer "Ready for handoff" status text) ---
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth     = 16.33
$dede.Columns.Item(3).ColumnWidth     = 16.33
